# Reorder the 9 time-slot rows (rows 2-10) into chronological order based on
# the time-slot label in column B, keeping each row's A (order index) and
# C:G (per-day availability) values attached to its own B label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "C", "D", "E", "F", "G")

function Get-MinutesFromLabel($label) {
    $startText = ($label -split " - ")[0].Trim()
    # $startText looks like "8:00am" or "12:30pm"
    $ampm = $startText.Substring($startText.Length - 2, 2).ToLower()
    $timePart = $startText.Substring(0, $startText.Length - 2)
    $parts = $timePart -split ":"
    $hh = [int]$parts[0]
    $mm = [int]$parts[1]
    if ($ampm -eq "pm" -and $hh -ne 12) {
        $hh = $hh + 12
    }
    if ($ampm -eq "am" -and $hh -eq 12) {
        $hh = 0
    }
    return ($hh * 60 + $mm)
}

# Read the 9 existing data rows (rows 2-10) into memory first.
$records = @()
for ($r = 2; $r -le 10; $r++) {
    $rec = @{}
    foreach ($col in $cols) {
        $rec[$col] = $ws.Range("$col$r").Value2
    }
    $rec["SortKey"] = Get-MinutesFromLabel $rec["B"]
    $records += $rec
}

$sortedRecords = $records | Sort-Object -Property SortKey

# Write the chronologically sorted rows back into rows 2-10.
$r = 2
foreach ($rec in $sortedRecords) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rec[$col]
    }
    $r++
}
